# Update cryptos list: Price (column D) and Volume(1h) (column E)
# for rows 2-51, reflecting refreshed market data.
# Number format is forced to Text ("@") before assignment so that
# numeric-looking strings (e.g. "233.55") stay stored as text, matching
# the original inlineStr cell type used throughout this sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.760.23'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.085.86'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.55'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.638'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.24%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '58.06'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.392'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0781'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.34%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.17'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.393.15'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.07'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.777'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.34'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.078.80'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.743.96'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.10'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.95'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.04'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.60%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.38'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.94%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.72'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +8.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.76'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.134'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.51'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.38'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.48%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.68'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.90%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.48'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.67%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.39%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0235'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +9.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '101.07'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0969'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.87%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.35%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.451.34'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.63'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.11'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.05'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.20'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.30%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.276.97'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.49%  '
